$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 200 ("Fruta / hortaliza, semanal").
# Insert a blank row at 200, pushing the former rows 200..299 down to 201..300,
# then populate the newly inserted row with the new record's data.
$ws.Rows("200:200").Insert()

$ws.Range("A200").Value = 4
$ws.Range("B200").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C200").Value = "Los Lagos"
$ws.Range("D200").Value = 45029
$ws.Range("E200").Value = 10
$ws.Range("F200").Value = 100112009
$ws.Range("G200").Value = "Acelga"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 25
$ws.Range("K200").Value = 9000
$ws.Range("L200").Value = 9000
$ws.Range("M200").Value = 9000
$ws.Range("N200").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O200").Value = "Región de La Araucanía"
$ws.Range("P200").Value = 750
$ws.Range("Q200").Value = 12
$ws.Range("R200").Value = "Hortaliza"
